# Update DefinedTerm and DefinedTermSet classes and regen project
#
# - Renames the existing three sheets to their new schema.org class names.
# - Inserts "alternateName" / "sameAs" / "url" columns (and re-orders
#   "description"/"name") on the Thing / Intangible / CreativeWork sheets.
# - Adds two brand-new sheets, DefinedTerm and DefinedTermSet, each with
#   their own header row.

$wb = $excel.ActiveWorkbook

# --- Rename the pre-existing sheets -----------------------------------
$wb.Worksheets.Item(1).Name = "Thing"
$wb.Worksheets.Item(2).Name = "Intangible"
$wb.Worksheets.Item(3).Name = "CreativeWork"

# --- Shared header row used by Thing / Intangible / CreativeWork -------
$commonHeaders = @("id", "alternateName", "description", "name", "sameAs", "url")

foreach ($sheetName in @("Thing", "Intangible", "CreativeWork")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $commonHeaders.Count; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $commonHeaders[$i]
    }
}

# --- New sheet: DefinedTerm ---------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$definedTerm = $wb.Worksheets.Add($null, $lastSheet)
$definedTerm.Name = "DefinedTerm"
$definedTermHeaders = @("inDefinedTermSet", "termCode", "id", "alternateName", "description", "name", "sameAs", "url")
for ($i = 0; $i -lt $definedTermHeaders.Count; $i++) {
    $definedTerm.Cells.Item(1, $i + 1).Value = $definedTermHeaders[$i]
}

# --- New sheet: DefinedTermSet ------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$definedTermSet = $wb.Worksheets.Add($null, $lastSheet)
$definedTermSet.Name = "DefinedTermSet"
$definedTermSetHeaders = @("hasDefinedTerm", "id", "alternateName", "description", "name", "sameAs", "url")
for ($i = 0; $i -lt $definedTermSetHeaders.Count; $i++) {
    $definedTermSet.Cells.Item(1, $i + 1).Value = $definedTermSetHeaders[$i]
}
